# Update LoopUnrollingFactor2 clock-enable power report with refreshed
# Vivado utilization numbers, and drop the two lowest-power rows (old
# rows 9 and 10) that no longer appear in the refreshed report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.00009638717165216804

$ws.Range("A3").Value = 0.00003548925451468676
$ws.Range("C3").Value = 36.417911529541016
$ws.Range("D3").Value = 18.208955764770508
$ws.Range("E3").Value = 39.0
$ws.Range("F3").Value = 12.0
$ws.Range("H3").Value = "FF LUT "

$ws.Range("A4").Value = 0.000020157314793323167
$ws.Range("B4").Value = "firConvolutionLoopUnrollingFactor2_IP/U0/ap_CS_fsm_state6"
$ws.Range("C4").Value = 36.417911529541016
$ws.Range("D4").Value = 18.208955764770508
$ws.Range("E4").Value = 19.0
$ws.Range("F4").Value = 6.0
$ws.Range("H4").Value = "DSP FF "

$ws.Range("A5").Value = 0.0000188095527846599
$ws.Range("C5").Value = 39.850746154785156
$ws.Range("D5").Value = 19.99140167236328
$ws.Range("E5").Value = 10.0
$ws.Range("F5").Value = 6.0

$ws.Range("A6").Value = 0.000014020895832800306
$ws.Range("B6").Value = "firConvolutionLoopUnrollingFactor2_IP/U0/ap_CS_fsm_state5"
$ws.Range("C6").Value = 36.417911529541016
$ws.Range("D6").Value = 18.208955764770508
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 5.0

$ws.Range("A7").Value = 0.0000063626866904087365
$ws.Range("B7").Value = "firConvolutionLoopUnrollingFactor2_IP/U0/i_2_reg_2800"
$ws.Range("C7").Value = 36.56716537475586
$ws.Range("D7").Value = 18.33510398864746
$ws.Range("E7").Value = 6.0
$ws.Range("F7").Value = 3.0

$ws.Range("A8").Value = 0.0000015474626025024918
$ws.Range("C8").Value = 21.492536544799805
$ws.Range("D8").Value = 8.95522403717041
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0

# Remove rows 9 and 10 (now both collapse to row 9 after each delete)
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()

